# Update for alt key.
# Insert two new rows (for the new "alt key view flag" / "ak" entries) right
# before the existing "widget style" rows (which were rows 23-24, 1-indexed
# including the header), shifting everything below down by two rows, and
# renumber the "No." column (A) to keep it sequential.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert two blank rows before row 23.
$ws.Rows.Item(23).Resize(2).Insert()

# Fill in the two new rows with the "alt key view flag" data, following the
# same pattern used by the other "*** key view flag" boolean pairs elsewhere
# in the sheet (Default "-" => Value "true", Default "y" => Value "false").
$ws.Range("A23").Value = 22
$ws.Range("B23").Value = "alt key view flag"
$ws.Range("C23").Value = "ak"
$ws.Range("D23").Value = "boolean"
$ws.Range("E23").Value = 1
$ws.Range("F23").Value = "-"
$ws.Range("G23").Value = "'true"

$ws.Range("A24").Value = 23
$ws.Range("B24").Value = "alt key view flag"
$ws.Range("C24").Value = "ak"
$ws.Range("D24").Value = "boolean"
$ws.Range("E24").Value = 1
$ws.Range("F24").Value = "y"
$ws.Range("G24").Value = "'false"

# Renumber the "No." column for every row from the old row 23 (now row 25)
# through the end of the table so the sequence stays contiguous.
$lastRow = $ws.UsedRange.Rows.Count
for ($r = 25; $r -le $lastRow; $r++) {
    $ws.Cells.Item($r, 1).Value = $r - 1
}
